$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed coin values scraped on Fri Nov 22 19:22:36 UTC 2024.
# Price/Volume columns (D/E) are forced to text format per-cell so that
# values such as "1.00", "0.998" or "99.400.39" are stored verbatim as
# text instead of being re-interpreted by Excel as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '99.400.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.69%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.275.24'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.35%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.44'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '620.72'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.41'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +17.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.397'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.98%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.952'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +18.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.273.95'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.10'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +9.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '99.011.53'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000246'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.867.28'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.44'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.268.89'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.41'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.21'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.25'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.72'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.25'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000199'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.60'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.63'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +26.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.86'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.424.11'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.22%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.25'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +10.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.63'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.55%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.469'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.15'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.25%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.83'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '487.95'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.62%  '
$ws.Range("B42").Value = 'MantraDAO'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.62'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.22'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.72%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.766'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.06'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.94'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '157.31'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.840'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.20'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +14.60%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.67'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.87%  '

Write-Host "Updated 93 cells."
